$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; this shifts existing rows 5-37 down to 6-38
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new record's data
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(5, 3).Value = "Los Lagos"
$ws.Cells.Item(5, 4).Value = 44635
$ws.Cells.Item(5, 5).Value = 10
$ws.Cells.Item(5, 6).Value = 100112030
$ws.Cells.Item(5, 7).Value = "Poroto granado"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 70
$ws.Cells.Item(5, 11).Value = 32000
$ws.Cells.Item(5, 12).Value = 32000
$ws.Cells.Item(5, 13).Value = 32000
$ws.Cells.Item(5, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(5, 15).Value = "Región Metropolitana"
$ws.Cells.Item(5, 16).Value = 1280
$ws.Cells.Item(5, 17).Value = 25
$ws.Cells.Item(5, 18).Value = "Hortaliza"
